{"js": "// The document has two tables: the first (\"Zeile 1/2, Spalte 1/2\") is a\n// duplicate demo table that should be removed, along with the now-useless\n// empty paragraph that sits directly above it. The second table (\"Column 1/2\"\n// header row) and everything else stays untouched.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst firstTable = tables.items[0];\nconst paraBeforeFirstTable = firstTable.getParagraphBefore();\nparaBeforeFirstTable.load(\"text\");\nawait context.sync();\n\n// Sanity guard: only delete the leading paragraph if it is indeed empty\n// (matches the `<w:p/>` that directly precedes the removed table).\nif (paraBeforeFirstTable.text === \"\") {\n  paraBeforeFirstTable.delete();\n}\nfirstTable.delete();\nawait context.sync();\n", "ps1": "# The document has two tables: the first one (\"Zeile 1/2, Spalte 1/2\") is a\n# duplicate demo table that should be removed, along with the now-useless\n# empty paragraph that sits directly above it. The second table (\"Column 1/2\"\n# header row) and everything else stays untouched.\n$d = $word.ActiveDocument\n\n$firstTable = $d.Tables(1)\n\n# Paragraph immediately preceding the first table (the empty `<w:p/>`).\n$paraBeforeFirstTable = $d.Range(0, $firstTable.Range.Start).Paragraphs.Last\n\n# Sanity guard: only delete the leading paragraph if it is indeed empty\n# (just the paragraph mark, matching the `<w:p/>` that directly precedes the\n# removed table).\nif ($paraBeforeFirstTable.Range.Text -eq \"`r\") {\n    $paraBeforeFirstTable.Range.Delete()\n}\n\n$firstTable.Delete()\n"}
